$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.801.29"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +5.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.643.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +5.80%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.16"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.76%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "195.01"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.636.05"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.78%  "

$ws.Range("E9").Value = "  -0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.186"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +8.61%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.680"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.13"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000305"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +10.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.00"
$ws.Range("D14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.233.53"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +5.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.45"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.39%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.649.90"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.851.74"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +5.78%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.63%  "

$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "489.31"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.29"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +13.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("E25").Value = "  +2.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.49"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.18"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +6.12%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.40"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.59"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.33%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.123"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +10.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.30"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "618.53"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +3.19%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "66.43"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.12%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.23"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.05%  "

$ws.Range("E37").Value = "  +11.00%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.412"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.36%  "

$ws.Range("E39").Value = "  +1.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.57"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.337.36"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.22"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +17.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.18"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +8.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +8.28%  "

$ws.Range("E46").Value = "  +6.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.64"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +10.98%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.140"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.58%  "

$ws.Range("E50").Value = "  +1.56%  "

$ws.Range("E51").Value = "  -0.17%  "
